$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = '2024-02-10 00:29:47'
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.001
$ws.Range("J14").Value = 0.05
$ws.Range("K14").Value = 0.003
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 500
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 5
$ws.Range("P14").Value = 2
$ws.Range("Q14").Value = 500
$ws.Range("R14").Value = 6
$ws.Range("S14").Value = 6
$ws.Range("T14").Value = 150
$ws.Range("U14").Value = 0.9166666666666666
$ws.Range("V14").Value = 'Data/bombay1.xlsx'
$ws.Range("W14").Value = 229500

$ws.Range("A15").Value = '2024-02-10 18:10:16'
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.001
$ws.Range("J15").Value = 0.05
$ws.Range("K15").Value = 0.003
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 5
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 500
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("T15").Value = 150
$ws.Range("U15").Value = 0.9
$ws.Range("V15").Value = 'Data/bombay1.xlsx'
$ws.Range("W15").Value = 255000

$ws.Range("A16").Value = '2024-02-12 08:25:39'
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0.001
$ws.Range("J16").Value = 0.05
$ws.Range("K16").Value = 0.003
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = 500
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 5
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 1000
$ws.Range("R16").Value = 30
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 100
$ws.Range("U16").Value = 1
$ws.Range("V16").Value = 'Data/bombay1.xlsx'
$ws.Range("W16").Value = 1581000

$ws.Range("A17").Value = '2024-02-12 08:38:06'
$ws.Range("B17").Value = 18
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0.001
$ws.Range("J17").Value = 0.05
$ws.Range("K17").Value = 0.003
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 500
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 5
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 200
$ws.Range("R17").Value = 20
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 100
$ws.Range("U17").Value = 0.8888888888888888
$ws.Range("V17").Value = 'Data/bombayauto.xlsx'
$ws.Range("W17").Value = -21600

$ws.Range("A18").Value = '2024-02-12 11:36:44'
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0.001
$ws.Range("J18").Value = 0.05
$ws.Range("K18").Value = 0.003
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = 500
$ws.Range("N18").Value = 10
$ws.Range("O18").Value = 5
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 1000
$ws.Range("R18").Value = 12
$ws.Range("S18").Value = 6
$ws.Range("T18").Value = 200
$ws.Range("U18").Value = 0.9166666666666666
$ws.Range("V18").Value = 'Data/bombayauto.xlsx'
$ws.Range("W18").Value = 148000

$ws.Range("A19").Value = '2024-02-12 12:03:08'
$ws.Range("B19").Value = 19
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0.001
$ws.Range("J19").Value = 0.05
$ws.Range("K19").Value = 0.003
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = 500
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 5
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 1000
$ws.Range("R19").Value = 8
$ws.Range("S19").Value = 6
$ws.Range("T19").Value = 200
$ws.Range("U19").Value = 0.7894736842105263
$ws.Range("V19").Value = 'Data/bombayauto.xlsx'
$ws.Range("W19").Value = 510000

$ws.Range("A20").Value = '2024-02-12 16:48:11'
$ws.Range("B20").Value = 32
$ws.Range("C20").Value = 26
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 11
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0.001
$ws.Range("J20").Value = 0.05
$ws.Range("K20").Value = 0.003
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 500
$ws.Range("N20").Value = 10
$ws.Range("O20").Value = 5
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 500
$ws.Range("R20").Value = 9
$ws.Range("S20").Value = 6
$ws.Range("T20").Value = 150
$ws.Range("U20").Value = 0.8125
$ws.Range("V20").Value = 'Data/bombay1.xlsx'
$ws.Range("W20").Value = 52500

$ws.Range("A21").Value = '2024-02-13 14:04:23'
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0.001
$ws.Range("J21").Value = 0.05
$ws.Range("K21").Value = 0.003
$ws.Range("L21").Value = 100
$ws.Range("M21").Value = 500
$ws.Range("N21").Value = 10
$ws.Range("O21").Value = 5
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 1000
$ws.Range("R21").Value = 4
$ws.Range("S21").Value = 7
$ws.Range("T21").Value = 150
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 'Data/bombayauto.xlsx'
$ws.Range("W21").Value = -111000
